$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update environment host/url for the "pre prod" row (row 3) to the new
# "i-preproduccion" host used by the regression (R33) data.
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"

# Keep the existing hyperlink target (on B3) in sync with the new URL,
# updating it in place instead of creating a new hyperlink relationship.
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$3') {
        $hl.Address = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
    }
}

# Update the Documento and NumeroCalle values used for the account creation data.
$ws.Range("G3").Value = 30629625131
$ws.Range("M3").Value = 308
